$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    for ($col = 2; $col -le 29; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

Swap-Rows 17 18
Swap-Rows 58 59
Swap-Rows 73 74
Swap-Rows 89 90
Swap-Rows 91 92
Swap-Rows 103 104
Swap-Rows 151 152
